$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three user story descriptions with expanded client answer text
$ws.Range("C5").Value = "Whenever I enter a time and pace, the app tells me how far I could run. In miles."

$ws.Range("C4").Value = "Whenever I enter a time and event (5K, 5-mile, 10K, half-marathon), the app tells me what my pace should be (during prep) or was (if completed). Pace must be displayed as minutes, seconds, and tenths of seconds per mile, e.g., 7:03.6" + [char]10 + "Pace must have an upper limit of 20 min/mile. Only one decimal digit."

$ws.Range("C6").Value = "Whenever I enter an event (5K, 5-mile, 10K, half-marathon), and pace, the app tells me how much time it would take to run that distance.  Time must be displayed as hours, minutes, seconds, and tenths of seconds, e.g., 1:43:54.7" + [char]10 + "Time has an upper limit of 4 hours.  Only one decimal digit.  Only whole values for hours are allowed."

# Adjust row heights to match expanded content
$ws.Rows("4").RowHeight = 126
$ws.Rows("6").RowHeight = 120

# Update the active selection
[void]$ws.Range("C7").Select()
